$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell J1: shared string "2024", same style as the other
# --- header cells (s="1"). Format the cell as Text first so the numeric-
# --- looking literal is stored as a shared string (t="s") rather than
# --- being auto-parsed into a number, then copy the visual formatting
# --- (border/font/alignment) from I1 so J1 matches the rest of the header
# --- row.
$ws.Range("J1").NumberFormat = "@"
$ws.Range("J1").Value = "2024"
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# --- New data column J (2024 figures) for rows 2-35 ---
$ws.Range("J2").Value = 4045.6937499999999
$ws.Range("J3").Value = 3814.9729166666671
$ws.Range("J4").Value = 271.91979166666664
$ws.Range("J5").Value = 268.45937499999997
$ws.Range("J6").Value = 201.82291666666666
$ws.Range("J7").Value = 13.505208333333334
$ws.Range("J8").Value = 2.1802083333333333
$ws.Range("J9").Value = 4.7437499999999995
$ws.Range("J10").Value = 395.75104166666665
$ws.Range("J11").Value = 1114.7437499999999
$ws.Range("J12").Value = 113.25625000000001
$ws.Range("J13").Value = 1.6156249999999999
$ws.Range("J14").Value = 264.09270833333335
$ws.Range("J15").Value = 3.4010416666666665
$ws.Range("J16").Value = 6.2
$ws.Range("J17").Value = 5.8
$ws.Range("J18").Value = 3.6
$ws.Range("J19").Value = 6.7
$ws.Range("J20").Value = 7
$ws.Range("J21").Value = 2.2000000000000002
$ws.Range("J22").Value = 94.3
$ws.Range("J23").Value = 5.8
$ws.Range("J24").Value = 16.100000000000001
$ws.Range("J25").Value = 11
$ws.Range("J26").Value = 65
$ws.Range("J27").Value = 5.3
$ws.Range("J28").Value = 5.0999999999999996
$ws.Range("J29").Value = 2.8
$ws.Range("J30").Value = 702.2
$ws.Range("J31").Value = 19.7
$ws.Range("J32").Value = 19.8
$ws.Range("J33").Value = 81.7
$ws.Range("J34").Value = 10.9
$ws.Range("J35").Value = 10.4

# --- Match the saved selection/active cell ---
$ws.Range("L3").Select()
